# ACG_Common_Workbook.xlsx edit
# commit message: "change the values for product, location, user, role"
#
# Rights_test_case : Admin22/Creator22      -> Admin23/Creator23
# Users             : ACG3300..ACG3307      -> ACG3400..ACG3407 (userid + email)
# partner           : location15..location45 -> location16..location46
#                     364450..364453         -> 364460..364463
# Product           : 515391..515394         -> 515401..515404
#                     Product20..Product23    -> Product30..Product33

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Rights_test_case sheet
# ---------------------------------------------------------------
$wsRights = $wb.Worksheets.Item("Rights_test_case")
$wsRights.Range("C8").Value = "Admin23"
$wsRights.Range("C9").Value = "Creator23"

# ---------------------------------------------------------------
# Users sheet
# ---------------------------------------------------------------
$wsUsers = $wb.Worksheets.Item("Users")
for ($i = 0; $i -lt 8; $i++) {
    $row = 2 + $i
    $oldNum = 3300 + $i
    $newNum = 3400 + $i
    $newId = "ACG$newNum"
    $newEmail = "ACG$newNum@gmail.com"

    $wsUsers.Range("A$row").Value = $newId
    $wsUsers.Range("B$row").Value = $newEmail
    $wsUsers.Range("P$row").Value = $newId
}

# ---------------------------------------------------------------
# partner sheet
# ---------------------------------------------------------------
$wsPartner = $wb.Worksheets.Item("partner")
$partnerLocations = @("location16", "location26", "location36", "location46")
for ($i = 0; $i -lt 4; $i++) {
    $row = 2 + $i
    $newLoc = $partnerLocations[$i]
    $newIdent = 364460 + $i

    $wsPartner.Range("A$row").Value = $newLoc
    $wsPartner.Range("B$row").Value = $newIdent
}

# ---------------------------------------------------------------
# Product sheet
# ---------------------------------------------------------------
$wsProduct = $wb.Worksheets.Item("Product")
for ($i = 0; $i -lt 4; $i++) {
    $row = 2 + $i
    $newNum = 515401 + $i
    $newName = "Product" + (30 + $i)

    $wsProduct.Range("A$row").Value = $newNum
    $wsProduct.Range("B$row").Value = $newName
    $wsProduct.Range("E$row").Value = $newNum
    $wsProduct.Range("G$row").Value = $newNum
}

# ---------------------------------------------------------------
# View / selection state
# ---------------------------------------------------------------
[void]$wsRights.Activate()
[void]$wsRights.Range("C9").Select()

[void]$wsUsers.Activate()
[void]$wsUsers.Range("D15").Select()

[void]$wsPartner.Activate()
[void]$wsPartner.Range("E12").Select()

[void]$wsProduct.Activate()
[void]$wsProduct.Range("G7").Select()
